# Daily attendance processing - 2025-11-29 18:31:50
# Reorders the "Recorded By" (column G) author lists on the
# "Session Analysis Results" sheet:
#   "System, dnasr281@gmail.com"            -> "dnasr281@gmail.com, System"
#   "backup@backdoor.com, System, system"   -> "system, backup@backdoor.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Text

    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
    elseif ($val -eq "backup@backdoor.com, System, system") {
        $cell.Value = "system, backup@backdoor.com, System"
    }
}
